$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# 1. Remove the "Meta description" paragraph (currently paragraph #2,
#    right after the title heading "Play Diego Wild for Free...").
$metaPara = $d.Paragraphs.Item(2)
$null = $metaPara.Range.Delete()

# 2. Insert a new bold paragraph ("Play Diego Wild for Free: See Gameplay,
#    Payouts & Bonuses") right before the final (DALLE-prompt) paragraph.
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$null = $insertPoint.InsertParagraphBefore()

# The paragraph that was just created now sits at $lastIndex (the old last
# paragraph shifted down to $lastIndex + 1). Its range currently contains
# only the paragraph mark; replace that with the desired run markup.
$newPara = $d.Paragraphs.Item($lastIndex)
$newXml = "<w:p $wNs><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Diego Wild for Free: See Gameplay, Payouts &amp; Bonuses</w:t></w:r></w:p>"
$null = $newPara.Range.InsertXML($newXml)

# 3. Replace the text of the final (DALLE-prompt) paragraph with the new
#    meta-description text, keeping its italic formatting and run layout.
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalRange = $finalPara.Range
$finalTrimmed = $d.Range($finalRange.Start, $finalRange.End - 1)
$finalTrimmed.Text = "Find out more about Diego Wild, the slot game set in the Amazon forest. Play for free or real money and explore ancient Aztec temples."
